# Apply the "1st changes of mifos to finflux" edit:
# - Insert a new blank column at N on the "Repayment Schedule" sheet,
#   shifting the old N/O/P columns to O/P/Q.
# - Make "Repayment Schedule" the active sheet (instead of "NewLoanInput"),
#   and set the active selection on that sheet to S5.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before column N (pushes N,O,P -> O,P,Q)
$wsSchedule.Columns("N").Insert()

# Activate the Repayment Schedule sheet and set selection to S5
$wsSchedule.Activate()
$wsSchedule.Range("S5").Select()
